$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 9: CCMastercard ----
$ws.Range("A9").Value = "CCMastercard"
$ws.Range("O9").Value = "'5555555555554444"
$null = $ws.Range("P3").Copy()
$null = $ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = "'04/26"
$ws.Range("Q9").Value = 123

# ---- Row 10: CCAmexcard ----
$ws.Range("A10").Value = "CCAmexcard"
$ws.Range("O10").Value = "'378282246310005"
$null = $ws.Range("P3").Copy()
$null = $ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = "'04/26"
$ws.Range("Q10").Value = 1234

# ---- Row 11: CCDiscovercard ----
$ws.Range("A11").Value = "CCDiscovercard"
$ws.Range("O11").Value = "'6011111111111117"
$null = $ws.Range("P3").Copy()
$null = $ws.Range("P11").PasteSpecial(-4122)
$ws.Range("P11").Value = "'04/26"
$ws.Range("Q11").Value = 123

$excel.CutCopyMode = 0

# ---- View / window state ----
# Move the visible top-left cell back toward the front of the sheet and
# update the active selection, matching the refreshed sheetView.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$null = $ws.Range("N14").Select()

# Workbook window placement/size (best-effort; mirrors workbookView change).
$win.Left = 1335
$win.Top = 1500
$win.Width = 21600
$win.Height = 9675
